$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 8 new rows before row 5 (shift existing rows 5-26 down to 13-34)
$ws.Rows("5:12").Insert(-4121) | Out-Null

# 2. Fill the newly inserted rows (5-12) with new benchmark data
$ws.Cells.Item(5,1).Value = "StandardScaler"
$ws.Cells.Item(5,2).Value = 3
$ws.Cells.Item(5,3).Value = "[128, 128, 1]"
$ws.Cells.Item(5,4).Value = "['relu', 'relu', 'linear']"
$ws.Cells.Item(5,5).Value = "Adam"
$ws.Cells.Item(5,6).Value = 15
$ws.Cells.Item(5,7).Value = 10000
$ws.Cells.Item(5,8).Value = 0.3605457623799642
$ws.Cells.Item(5,9).Value = 5.408186435699463
$ws.Cells.Item(5,10).Value = 11.92198085784912

$ws.Cells.Item(6,1).Value = "StandardScaler"
$ws.Cells.Item(6,2).Value = 3
$ws.Cells.Item(6,3).Value = "[128, 128, 1]"
$ws.Cells.Item(6,4).Value = "['relu', 'relu', 'linear']"
$ws.Cells.Item(6,5).Value = "Adam"
$ws.Cells.Item(6,6).Value = 154
$ws.Cells.Item(6,7).Value = 10000
$ws.Cells.Item(6,8).Value = 0.391872331693575
$ws.Cells.Item(6,9).Value = 60.34833908081055
$ws.Cells.Item(6,10).Value = 4.974963188171387

$ws.Cells.Item(7,1).Value = "StandardScaler"
$ws.Cells.Item(7,2).Value = 3
$ws.Cells.Item(7,3).Value = "[128, 128, 1]"
$ws.Cells.Item(7,4).Value = "['relu', 'relu', 'linear']"
$ws.Cells.Item(7,5).Value = "Adam"
$ws.Cells.Item(7,6).Value = 157
$ws.Cells.Item(7,7).Value = 10000
$ws.Cells.Item(7,8).Value = 0.4601679030497363
$ws.Cells.Item(7,9).Value = 72.24636077880859
$ws.Cells.Item(7,10).Value = 2.854086399078369

$ws.Cells.Item(8,1).Value = "StandardScaler"
$ws.Cells.Item(8,2).Value = 3
$ws.Cells.Item(8,3).Value = "[128, 128, 1]"
$ws.Cells.Item(8,4).Value = "['relu', 'relu', 'linear']"
$ws.Cells.Item(8,5).Value = "Adam"
$ws.Cells.Item(8,6).Value = 132
$ws.Cells.Item(8,7).Value = 10000
$ws.Cells.Item(8,8).Value = 0.3442076715556058
$ws.Cells.Item(8,9).Value = 45.43541264533997
$ws.Cells.Item(8,10).Value = 0.3319054841995239

$ws.Cells.Item(9,1).Value = "StandardScaler"
$ws.Cells.Item(9,2).Value = 3
$ws.Cells.Item(9,3).Value = "[128, 128, 1]"
$ws.Cells.Item(9,4).Value = "['relu', 'relu', 'linear']"
$ws.Cells.Item(9,5).Value = "Adam"
$ws.Cells.Item(9,6).Value = 85
$ws.Cells.Item(9,7).Value = 10000
$ws.Cells.Item(9,8).Value = 0.9589118817273308
$ws.Cells.Item(9,9).Value = 81.50750994682312
$ws.Cells.Item(9,10).Value = 0.004373230971395969

$ws.Cells.Item(10,1).Value = "StandardScaler"
$ws.Cells.Item(10,2).Value = 3
$ws.Cells.Item(10,3).Value = "[128, 128, 1]"
$ws.Cells.Item(10,4).Value = "['relu', 'relu', 'linear']"
$ws.Cells.Item(10,5).Value = "Adam"
$ws.Cells.Item(10,6).Value = 69
$ws.Cells.Item(10,7).Value = 10000
$ws.Cells.Item(10,8).Value = 0.7403726232224617
$ws.Cells.Item(10,9).Value = 51.08571100234985
$ws.Cells.Item(10,10).Value = 35.91139602661133

$ws.Cells.Item(11,1).Value = "StandardScaler"
$ws.Cells.Item(11,2).Value = 3
$ws.Cells.Item(11,3).Value = "[128, 128, 1]"
$ws.Cells.Item(11,4).Value = "['relu', 'relu', 'linear']"
$ws.Cells.Item(11,5).Value = "Adam"
$ws.Cells.Item(11,6).Value = 71
$ws.Cells.Item(11,7).Value = 10000
$ws.Cells.Item(11,8).Value = 0.4907066553411349
$ws.Cells.Item(11,9).Value = 34.84017252922058
$ws.Cells.Item(11,10).Value = 0.2965547144412994

$ws.Cells.Item(12,1).Value = "StandardScaler"
$ws.Cells.Item(12,2).Value = 3
$ws.Cells.Item(12,3).Value = "[128, 128, 1]"
$ws.Cells.Item(12,4).Value = "['relu', 'relu', 'linear']"
$ws.Cells.Item(12,5).Value = "Adam"
$ws.Cells.Item(12,6).Value = 145
$ws.Cells.Item(12,7).Value = 10000
$ws.Cells.Item(12,8).Value = 0.4851904753980966
$ws.Cells.Item(12,9).Value = 70.352618932724
$ws.Cells.Item(12,10).Value = 0.00001613486165297218

# 3. Add new header columns K1/L1 (copy formatting of existing header cell, then set values)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("K1:L1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(1,11).Value = "Unnamed: 10"
$ws.Cells.Item(1,12).Value = "Unnamed: 11"

# 4. Add the note in L12
$ws.Cells.Item(12,12).Value = "pozn F0 = 1N"
